# Natmi following Dr Hou advice
# Ligand/Receptor-expressing cell counts (cols E & K) go from 1 -> 3 for every
# data row, which ripples through every NATMI-derived statistic in that row
# (average/total expression, derived-specificity, and edge-weight columns).
# All values below are the recomputed NATMI outputs for the new cell counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.312821333333335
$ws.Range("H2").Value = 18.938464
$ws.Range("I2").Value = 0.3104630857074662
$ws.Range("J2").Value = 0.3104630857074661
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.816292333333331
$ws.Range("N2").Value = 29.448877
$ws.Range("O2").Value = 0.3637676450248369
$ws.Range("P2").Value = 0.3637676450248368
$ws.Range("Q2").Value = 61.96849965610311
$ws.Range("R2").Value = 557.716496904928
$ws.Range("S2").Value = 0.1129364255549491
$ws.Range("T2").Value = 0.112936425554949

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.312821333333335
$ws.Range("H3").Value = 18.938464
$ws.Range("I3").Value = 0.3104630857074662
$ws.Range("J3").Value = 0.3104630857074661
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 11.28584166666666
$ws.Range("N3").Value = 33.857525
$ws.Range("O3").Value = 0.4182255281116336
$ws.Range("P3").Value = 0.4182255281116335
$ws.Range("Q3").Value = 71.24550203795556
$ws.Range("R3").Value = 641.2095183416001
$ws.Range("S3").Value = 0.1298435879791724
$ws.Range("T3").Value = 0.1298435879791724

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.312821333333335
$ws.Range("H4").Value = 18.938464
$ws.Range("I4").Value = 0.3104630857074662
$ws.Range("J4").Value = 0.3104630857074661
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.882927666666667
$ws.Range("N4").Value = 17.648783
$ws.Range("O4").Value = 0.2180068268635295
$ws.Range("P4").Value = 0.2180068268635295
$ws.Range("Q4").Value = 37.13787127659023
$ws.Range("R4").Value = 334.2408414893121
$ws.Range("S4").Value = 0.06768307217334471
$ws.Range("T4").Value = 0.06768307217334468

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.354969666666667
$ws.Range("H5").Value = 28.064909
$ws.Range("I5").Value = 0.4600752335690602
$ws.Range("J5").Value = 0.4600752335690601
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.816292333333331
$ws.Range("N5").Value = 29.448877
$ws.Range("O5").Value = 0.3637676450248369
$ws.Range("P5").Value = 0.3637676450248368
$ws.Range("Q5").Value = 91.83111701746587
$ws.Range("R5").Value = 826.4800531571929
$ws.Range("S5").Value = 0.1673604842496688
$ws.Range("T5").Value = 0.1673604842496688

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.354969666666667
$ws.Range("H6").Value = 28.064909
$ws.Range("I6").Value = 0.4600752335690602
$ws.Range("J6").Value = 0.4600752335690601
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 11.28584166666666
$ws.Range("N6").Value = 33.857525
$ws.Range("O6").Value = 0.4182255281116336
$ws.Range("P6").Value = 0.4182255281116335
$ws.Range("Q6").Value = 105.5787064544694
$ws.Range("R6").Value = 950.2083580902249
$ws.Range("S6").Value = 0.1924152075305034
$ws.Range("T6").Value = 0.1924152075305033

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.354969666666667
$ws.Range("H7").Value = 28.064909
$ws.Range("I7").Value = 0.4600752335690602
$ws.Range("J7").Value = 0.4600752335690601
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.882927666666667
$ws.Range("N7").Value = 17.648783
$ws.Range("O7").Value = 0.2180068268635295
$ws.Range("P7").Value = 0.2180068268635295
$ws.Range("Q7").Value = 55.03460987286078
$ws.Range("R7").Value = 495.311488855747
$ws.Range("S7").Value = 0.100299541788888
$ws.Range("T7").Value = 0.100299541788888

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.665774
$ws.Range("H8").Value = 13.997322
$ws.Range("I8").Value = 0.2294616807234737
$ws.Range("J8").Value = 0.2294616807234737
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.816292333333331
$ws.Range("N8").Value = 29.448877
$ws.Range("O8").Value = 0.3637676450248369
$ws.Range("P8").Value = 0.3637676450248368
$ws.Range("Q8").Value = 45.80060154526599
$ws.Range("R8").Value = 412.205413907394
$ws.Range("S8").Value = 0.08347073522021904
$ws.Range("T8").Value = 0.08347073522021903

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.665774
$ws.Range("H9").Value = 13.997322
$ws.Range("I9").Value = 0.2294616807234737
$ws.Range("J9").Value = 0.2294616807234737
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.28584166666666
$ws.Range("N9").Value = 33.857525
$ws.Range("O9").Value = 0.4182255281116336
$ws.Range("P9").Value = 0.4182255281116335
$ws.Range("Q9").Value = 52.65718661644999
$ws.Range("R9").Value = 473.9146795480499
$ws.Range("S9").Value = 0.09596673260195786
$ws.Range("T9").Value = 0.09596673260195783

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.665774
$ws.Range("H10").Value = 13.997322
$ws.Range("I10").Value = 0.2294616807234737
$ws.Range("J10").Value = 0.2294616807234737
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.882927666666667
$ws.Range("N10").Value = 17.648783
$ws.Range("O10").Value = 0.2180068268635295
$ws.Range("P10").Value = 0.2180068268635295
$ws.Range("Q10").Value = 27.448410951014
$ws.Range("R10").Value = 247.035698559126
$ws.Range("S10").Value = 0.05002421290129683
$ws.Range("T10").Value = 0.05002421290129682
